$wb = $excel.ActiveWorkbook

# Row 36 data for each of the 4 worksheets, taken from the target diff.
$rows = @(
    @{
        Sheet = "DE_LFT_#1"
        A = 45822.43604166667
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x70"
        E = "0x14"
        F = 380
        G = 759863127514710945038336.0
        H = 368
        I = 14
    },
    @{
        Sheet = "DE_LFT_#2"
        A = 45822.43604166667
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x70"
        E = "0xe"
        F = 380
        G = 568432987514711010443264.0
        H = 368
        I = 14
    },
    @{
        Sheet = "DE_PLT_#1"
        A = 45822.43604166667
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x80"
        E = "0x7"
        F = 130
        G = 568631262647113970876416.0
        H = 128
        I = 7
    },
    @{
        Sheet = "DE_PLT_#2"
        A = 45822.43604166667
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x80"
        E = "0x3"
        F = 130
        G = 985046333984776009023488.0
        H = 128
        I = 3
    }
)

foreach ($row in $rows) {
    $ws = $wb.Worksheets.Item($row.Sheet)

    $cellA = $ws.Range("A36")
    $cellA.Value = $row.A
    $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Range("B36").Value = $row.B
    $ws.Range("C36").Value = $row.C
    $ws.Range("D36").Value = $row.D
    $ws.Range("E36").Value = $row.E
    $ws.Range("F36").Value = $row.F
    $ws.Range("G36").Value = $row.G
    $ws.Range("H36").Value = $row.H
    $ws.Range("I36").Value = $row.I
}
